$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 (two store rows no longer needed)
$ws.Rows("3:4").Delete()

# Update remaining data row (row 2) with new values
# Keep A2 as text (store "name") like the original data, not a number.
# Temporarily force Text format so Excel doesn't coerce the numeric-looking
# string into a number, then clear the format override so the cell ends up
# with no explicit style (matching the rest of column A).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "13254"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 57
$ws.Range("D2").Value = 57100
$ws.Range("E2").Value = 37540
$ws.Range("F2").Value = -19560
$ws.Range("G2").Value = -0.3425569176882662
